# Commit: "subi iamgenes de bremen, agregue cta 60109"
# This workbook is a single-sheet client/permissions list (CUENTA, NOMBRE,
# CATEGORIAS, LISTA_PRECIOS) sorted by account number. We need to insert a
# new client account 60109 "CERAMICA SANTA MARTA" in its correct sorted
# position (right after account 60103, i.e. it becomes new row 235, and
# every following row shifts down by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$insertRow = 235

# 1) Insert a brand new blank row at 235, pushing 235..323 down to 236..324.
$ws.Rows.Item($insertRow).Insert()

# 2) The freshly inserted row has no formatting of its own; pull the exact
#    per-column formatting used throughout the table from neighboring rows
#    that already use the right style for each column:
#      - column C (CATEGORIAS) always uses the same style as any other row
#        in this block (style copied from the row directly above, 234)
#      - column D (LISTA_PRECIOS) needs the style used for value "F"; the
#        row that used to be 305 (70114 / FERRETERIA FERRECAS) is now row
#        306 after the insert above, and already carries that exact style.
$ws.Range("C234").Copy()
$ws.Range("C" + $insertRow).PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("D306").Copy()
$ws.Range("D" + $insertRow).PasteSpecial(-4122)   # xlPasteFormats

# Restore the standard row height / custom-height flag used by every data
# row in the table (new rows default to "no explicit height").
$ws.Rows.Item($insertRow).RowHeight = 15.75

# 3) Fill in the new account's data.
$ws.Range("A" + $insertRow).Value = 60109
$ws.Range("B" + $insertRow).Value = "CERAMICA SANTA MARTA"
$ws.Range("C" + $insertRow).Value = $ws.Range("C234").Value()
$ws.Range("D" + $insertRow).Value = "F"

# 4) The used range grew from A1:D323 to A1:D324 - refresh the AutoFilter
#    so its stored range (and the underlying defined name) matches.
$lastRow = 324
$ws.AutoFilterMode = $false
$ws.Range("A1:D" + $lastRow).AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Hoja 1!_FilterDatabase") {
        $n.RefersTo = "='Hoja 1'!`$A`$1:`$D`$" + $lastRow
    }
}

# 5) Match the author's final cursor/selection position in the sheet.
$ws.Range("B236").Select()
